$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = "[0.02138871637188977, 14.741770618790266]"
$ws.Range("N2").Value = [double]"0.04936562130411359"
$ws.Range("O2").Value = [double]"0.04936562130411359"
$ws.Range("Q2").Value = "[-2.792526803190927, 0.20126319302276974]"
$ws.Range("R2").Value = [double]"0.08811085872940994"
$ws.Range("S2").Value = [double]"0.08811085872940994"
$ws.Range("U2").Value = "[5.359680172348991, 13.606552239739901]"
$ws.Range("V2").Value = [double]"3.09769174577923e-05"
$ws.Range("W2").Value = [double]"3.09769174577923e-05"
$ws.Range("Y2").Value = [double]"-0.7332132132132179"
$ws.Range("Z2").Value = [double]"10.1733333333334"
$ws.Range("M3").Value = "[-0.7913786129725571, 15.900305250188381]"
$ws.Range("N3").Value = [double]"0.07493239867976409"
$ws.Range("O3").Value = [double]"0.07493239867976409"
$ws.Range("Q3").Value = "[-2.7422110049352346, 0.779894872963232]"
$ws.Range("R3").Value = [double]"0.2677534861523911"
$ws.Range("S3").Value = [double]"0.2677534861523911"
$ws.Range("U3").Value = "[6.377976570602483, 14.976227839566755]"
$ws.Range("V3").Value = [double]"9.134068148952679e-06"
$ws.Range("W3").Value = [double]"9.134068148952679e-06"
$ws.Range("Y3").Value = [double]"-2.84120120120122"
$ws.Range("Z3").Value = [double]"9.990030030030088"
$ws.Range("M4").Value = "[-0.7421529735876522, 14.1431708954099]"
$ws.Range("N4").Value = [double]"0.07646378508925489"
$ws.Range("O4").Value = [double]"0.07646378508925489"
$ws.Range("Q4").Value = "[-3.849158566560467, 2.213895123250465]"
$ws.Range("R4").Value = [double]"0.5896572595687504"
$ws.Range("S4").Value = [double]"0.5896572595687504"
$ws.Range("U4").Value = "[4.9872532715136275, 12.635061350556011]"
$ws.Range("V4").Value = [double]"3.00885158717179e-05"
$ws.Range("W4").Value = [double]"3.00885158717179e-05"
$ws.Range("Y4").Value = [double]"-8.065345345345394"
$ws.Range("Z4").Value = [double]"14.02270270270279"
$ws.Range("M5").Value = "[-2.3785081322639217, 13.773559576526342]"
$ws.Range("N5").Value = [double]"0.1622325243859912"
$ws.Range("O5").Value = [double]"0.1622325243859912"
$ws.Range("Q5").Value = "[-3.528395352680429, 2.3711319927995036]"
$ws.Range("R5").Value = [double]"0.6946426344892109"
$ws.Range("S5").Value = [double]"0.6946426344892109"
$ws.Range("U5").Value = "[5.107596134570322, 13.595733847229518]"
$ws.Range("V5").Value = [double]"5.816139283298938e-05"
$ws.Range("W5").Value = [double]"5.816139283298938e-05"
$ws.Range("Y5").Value = [double]"-8.638168168168219"
$ws.Range("Z5").Value = [double]"12.85414414414423"
$ws.Range("M6").Value = "[-1.6346358952862552, 12.351883259490325]"
$ws.Range("N6").Value = [double]"0.1297575860045488"
$ws.Range("O6").Value = [double]"0.1297575860045488"
$ws.Range("Q6").Value = "[-2.094395102393195, 1.6918687163476562]"
$ws.Range("R6").Value = [double]"0.8314177088751362"
$ws.Range("S6").Value = [double]"0.8314177088751362"
$ws.Range("U6").Value = "[4.179158905244932, 11.997039501901625]"
$ws.Range("V6").Value = [double]"0.0001379913979242531"
$ws.Range("W6").Value = [double]"0.0001379913979242531"
$ws.Range("Y6").Value = [double]"-6.163573573573612"
$ws.Range("Z6").Value = [double]"7.630000000000043"
$ws.Range("M7").Value = "[-1.743084331542331, 12.725689195362488]"
$ws.Range("N7").Value = [double]"0.1333094789002602"
$ws.Range("O7").Value = [double]"0.1333094789002602"
$ws.Range("Q7").Value = "[-2.8680005005744658, 3.182474239672544]"
$ws.Range("R7").Value = [double]"0.9170924588353047"
$ws.Range("S7").Value = [double]"0.9170924588353047"
$ws.Range("U7").Value = "[5.009618900869273, 13.540581650335636]"
$ws.Range("V7").Value = [double]"7.019709993327794e-05"
$ws.Range("W7").Value = [double]"7.019709993327794e-05"
$ws.Range("Y7").Value = [double]"11.29606606606613"
$ws.Range("Z7").Value = [double]"33.33828828828849"
$ws.Range("M8").Value = "[-0.7370032685369399, 13.519170752687634]"
$ws.Range("N8").Value = [double]"0.07763226809557411"
$ws.Range("O8").Value = [double]"0.07763226809557411"
$ws.Range("Q8").Value = "[-1.2830528555201557, 2.226474072814389]"
$ws.Range("R8").Value = [double]"0.5908848853857975"
$ws.Range("S8").Value = [double]"0.5908848853857975"
$ws.Range("U8").Value = "[4.77779218206036, 12.430517890165003]"
$ws.Range("V8").Value = [double]"4.332616402891176e-05"
$ws.Range("W8").Value = [double]"4.332616402891176e-05"
$ws.Range("Y8").Value = [double]"14.77882882882892"
$ws.Range("Z8").Value = [double]"27.5642342342344"
$ws.Range("M9").Value = "[-0.2909781330005483, 14.847183042466838]"
$ws.Range("N9").Value = [double]"0.05907940635982478"
$ws.Range("O9").Value = [double]"0.05907940635982478"
$ws.Range("Q9").Value = "[-0.7987632973091179, 2.7233425805893505]"
$ws.Range("R9").Value = [double]"0.2769391349113897"
$ws.Range("S9").Value = [double]"0.2769391349113897"
$ws.Range("U9").Value = "[5.460552922177817, 13.202908422666619]"
$ws.Range("V9").Value = [double]"1.487801919419773e-05"
$ws.Range("W9").Value = [double]"1.487801919419773e-05"
$ws.Range("Y9").Value = [double]"12.96870870870879"
$ws.Range("Z9").Value = [double]"25.7999399399401"
$ws.Range("M10").Value = "[0.23170266069744017, 16.275816038420263]"
$ws.Range("N10").Value = [double]"0.04399830587752795"
$ws.Range("O10").Value = [double]"0.04399830587752795"
$ws.Range("Q10").Value = "[-0.06289474781961513, 1.6604213424378491]"
$ws.Range("R10").Value = [double]"0.06841184520853005"
$ws.Range("S10").Value = [double]"0.06841184520853005"
$ws.Range("U10").Value = "[5.967295082912546, 14.342956872929687]"
$ws.Range("V10").Value = [double]"1.35211279048697e-05"
$ws.Range("W10").Value = [double]"1.35211279048697e-05"
$ws.Range("Y10").Value = [double]"16.84099099099109"
$ws.Range("Z10").Value = [double]"23.11912912912927"
$ws.Range("M11").Value = "[-0.07817963106551673, 14.51666805488333]"
$ws.Range("N11").Value = [double]"0.05240010354473057"
$ws.Range("O11").Value = [double]"0.05240010354473057"
$ws.Range("Q11").Value = "[-0.3459211130078863, 2.8239741771007356]"
$ws.Range("R11").Value = [double]"0.1223731210873809"
$ws.Range("S11").Value = [double]"0.1223731210873809"
$ws.Range("U11").Value = "[5.272436283113121, 13.40646718015501]"
$ws.Range("V11").Value = [double]"3.168162902400518e-05"
$ws.Range("W11").Value = [double]"3.168162902400518e-05"
$ws.Range("Y11").Value = [double]"12.882882882883"
$ws.Range("Z11").Value = [double]"24.68828828828853"
$ws.Range("M12").Value = "[-0.5618352151665675, 14.500691762856706]"
$ws.Range("N12").Value = [double]"0.06887504164256408"
$ws.Range("O12").Value = [double]"0.06887504164256408"
$ws.Range("Q12").Value = "[0.19497371824080734, 3.339711109221583]"
$ws.Range("R12").Value = [double]"0.02845205176787458"
$ws.Range("S12").Value = [double]"0.02845205176787458"
$ws.Range("U12").Value = "[5.372355350194303, 13.785065379757247]"
$ws.Range("V12").Value = [double]"3.593788084188532e-05"
$ws.Range("W12").Value = [double]"3.593788084188532e-05"
$ws.Range("Y12").Value = [double]"10.96216216216226"
$ws.Range("Z12").Value = [double]"22.67387387387409"
$ws.Range("M13").Value = "[-0.032050319587286324, 14.088667474425641]"
$ws.Range("N13").Value = [double]"0.05100510011806247"
$ws.Range("O13").Value = [double]"0.05100510011806247"
$ws.Range("Q13").Value = "[0.4968685077749626, 3.5661322013721986]"
$ws.Range("R13").Value = [double]"0.0106184541944796"
$ws.Range("S13").Value = [double]"0.0106184541944796"
$ws.Range("U13").Value = "[4.702548280961528, 12.334755065484343]"
$ws.Range("V13").Value = [double]"4.821098874807639e-05"
$ws.Range("W13").Value = [double]"4.821098874807639e-05"
$ws.Range("Y13").Value = [double]"10.11891891891901"
$ws.Range("Z13").Value = [double]"21.54954954954975"
$ws.Range("M14").Value = "[-0.4520051058652754, 14.387953855690673]"
$ws.Range("N14").Value = [double]"0.06501666629381142"
$ws.Range("O14").Value = [double]"0.06501666629381142"
$ws.Range("Q14").Value = "[0.5220264069028078, 3.5661322013721977]"
$ws.Range("R14").Value = [double]"0.00961317507407089"
$ws.Range("S14").Value = [double]"0.00961317507407089"
$ws.Range("U14").Value = "[4.723355302905166, 12.356085683189487]"
$ws.Range("V14").Value = [double]"4.655172700052113e-05"
$ws.Range("W14").Value = [double]"4.655172700052113e-05"
$ws.Range("Y14").Value = [double]"10.11891891891901"
$ws.Range("Z14").Value = [double]"21.45585585585606"
$ws.Range("M15").Value = "[-0.08073214936452189, 14.162961244812134]"
$ws.Range("N15").Value = [double]"0.05254257873027868"
$ws.Range("O15").Value = [double]"0.05254257873027868"
$ws.Range("Q15").Value = "[0.4339737599553466, 3.6290269491918146]"
$ws.Range("R15").Value = [double]"0.01385251725046843"
$ws.Range("S15").Value = [double]"0.01385251725046843"
$ws.Range("U15").Value = "[4.728990319940468, 12.389722381371477]"
$ws.Range("V15").Value = [double]"4.748844621027537e-05"
$ws.Range("W15").Value = [double]"4.748844621027537e-05"
$ws.Range("Y15").Value = [double]"9.884684684684771"
$ws.Range("Z15").Value = [double]"21.78378378378399"
$ws.Range("M16").Value = "[-0.9067824097166017, 15.027155952357266]"
$ws.Range("N16").Value = [double]"0.0810255029022926"
$ws.Range("O16").Value = [double]"0.0810255029022926"
$ws.Range("Q16").Value = "[-0.6541053773240009, 5.295737766411623]"
$ws.Range("R16").Value = [double]"0.1231287607806948"
$ws.Range("S16").Value = [double]"0.1231287607806948"
$ws.Range("U16").Value = "[5.470771630073827, 13.855476570222002]"
$ws.Range("V16").Value = [double]"2.994820800417131e-05"
$ws.Range("W16").Value = [double]"2.994820800417131e-05"
$ws.Range("Y16").Value = [double]"3.677477477477511"
$ws.Range("Z16").Value = [double]"25.83603603603628"
$ws.Range("M17").Value = "[-0.059080356088522024, 14.051225028447526]"
$ws.Range("N17").Value = [double]"0.05186770879032898"
$ws.Range("O17").Value = [double]"0.05186770879032898"
$ws.Range("Q17").Value = "[0.8868159442565786, 4.257974427387968]"
$ws.Range("R17").Value = [double]"0.003584744148775965"
$ws.Range("S17").Value = [double]"0.003584744148775965"
$ws.Range("U17").Value = "[4.812550917580567, 12.443924405528474]"
$ws.Range("V17").Value = [double]"3.98971112900437e-05"
$ws.Range("W17").Value = [double]"3.98971112900437e-05"
$ws.Range("Y17").Value = [double]"7.542342342342415"
$ws.Range("Z17").Value = [double]"20.09729729729749"
$ws.Range("M18").Value = "[-0.04196098858308517, 14.150154102250653]"
$ws.Range("N18").Value = [double]"0.05131270599062887"
$ws.Range("O18").Value = [double]"0.05131270599062887"
$ws.Range("Q18").Value = "[0.9371317425122707, 4.2076586291322755]"
$ws.Range("R18").Value = [double]"0.002754813068871265"
$ws.Range("S18").Value = [double]"0.002754813068871265"
$ws.Range("U18").Value = "[4.83666666322674, 12.463638708146942]"
$ws.Range("V18").Value = [double]"3.809607596250153e-05"
$ws.Range("W18").Value = [double]"3.809607596250153e-05"
$ws.Range("Y18").Value = [double]"7.729729729729804"
$ws.Range("Z18").Value = [double]"19.9099099099101"
$ws.Range("M19").Value = "[-0.696739308433262, 14.597225920451192]"
$ws.Range("N19").Value = [double]"0.07378854069652063"
$ws.Range("O19").Value = [double]"0.07378854069652063"
$ws.Range("Q19").Value = "[1.3019212798660411, 4.446658670846814]"
$ws.Range("R19").Value = [double]"0.000617973665811844"
$ws.Range("S19").Value = [double]"0.000617973665811844"
$ws.Range("U19").Value = "[5.4662288527698895, 13.828070818087303]"
$ws.Range("V19").Value = [double]"2.946286588567659e-05"
$ws.Range("W19").Value = [double]"2.946286588567659e-05"
$ws.Range("Y19").Value = [double]"6.839639639639707"
$ws.Range("Z19").Value = [double]"18.55135135135152"
$ws.Range("M20").Value = "[0.7552794845616244, 13.225858805582458]"
$ws.Range("N20").Value = [double]"0.0288382048518756"
$ws.Range("O20").Value = [double]"0.0288382048518756"
$ws.Range("Q20").Value = "[-4.39005339780916, -1.8239476867688489]"
$ws.Range("R20").Value = [double]"1.382611148570412e-05"
$ws.Range("S20").Value = [double]"1.382611148570412e-05"
$ws.Range("U20").Value = "[4.8243088147728805, 12.441284600843915]"
$ws.Range("V20").Value = [double]"3.849264086963622e-05"
$ws.Range("W20").Value = [double]"3.849264086963622e-05"
$ws.Range("Y20").Value = [double]"6.792792792792857"
$ws.Range("Z20").Value = [double]"16.3495495495497"
